# Update res_bus/vm_pu.xlsx results for Case_4_11 ("case with 380 kV done"):
# bus voltage magnitudes (rows 2-25, columns B-F and I-N) change because the
# slack-bus setpoint (column B) moved from 1.05 p.u. to 1.02 p.u., shifting
# every other bus's computed vm_pu accordingly. Column A (index), G (=1) and
# H (blank) are unaffected and left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.042376609325035
$ws.Cells.Item(2, 4).Value = 1.042710431440671
$ws.Cells.Item(2, 5).Value = 1.055695459013398
$ws.Cells.Item(2, 6).Value = 1.062600281309674
$ws.Cells.Item(2, 9).Value = 1.039716643391038
$ws.Cells.Item(2, 10).Value = 1.047452762640836
$ws.Cells.Item(2, 11).Value = 1.045486349988568
$ws.Cells.Item(2, 12).Value = 1.058435209496616
$ws.Cells.Item(2, 13).Value = 1.065321189121864
$ws.Cells.Item(2, 14).Value = 1.048940265708298

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.04332255882236
$ws.Cells.Item(3, 4).Value = 1.043399377259455
$ws.Cells.Item(3, 5).Value = 1.056662612344168
$ws.Cells.Item(3, 6).Value = 1.063708714444816
$ws.Cells.Item(3, 9).Value = 1.039934057007037
$ws.Cells.Item(3, 10).Value = 1.048045142904738
$ws.Cells.Item(3, 11).Value = 1.045986502198504
$ws.Cells.Item(3, 12).Value = 1.059215463282434
$ws.Cells.Item(3, 13).Value = 1.06624373195089
$ws.Cells.Item(3, 14).Value = 1.049533487220121

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.043935053734603
$ws.Cells.Item(4, 4).Value = 1.043845452508762
$ws.Cells.Item(4, 5).Value = 1.057289236617194
$ws.Cells.Item(4, 6).Value = 1.064427070793263
$ws.Cells.Item(4, 9).Value = 1.040073728857357
$ws.Cells.Item(4, 10).Value = 1.048428225592014
$ws.Cells.Item(4, 11).Value = 1.046309737435667
$ws.Cells.Item(4, 12).Value = 1.059720528971374
$ws.Cells.Item(4, 13).Value = 1.066841209977712
$ws.Cells.Item(4, 14).Value = 1.049917113928754

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.044192642123489
$ws.Cells.Item(5, 4).Value = 1.044033048672188
$ws.Cells.Item(5, 5).Value = 1.057552862702398
$ws.Cells.Item(5, 6).Value = 1.064729336485241
$ws.Cells.Item(5, 9).Value = 1.040132204833201
$ws.Cells.Item(5, 10).Value = 1.048589218452119
$ws.Cells.Item(5, 11).Value = 1.046445529766475
$ws.Cells.Item(5, 12).Value = 1.059932903052003
$ws.Cells.Item(5, 13).Value = 1.067092516040786
$ws.Cells.Item(5, 14).Value = 1.050078335417185

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.044235897926324
$ws.Cells.Item(6, 4).Value = 1.044064550736421
$ws.Cells.Item(6, 5).Value = 1.0575971379846
$ws.Cells.Item(6, 6).Value = 1.064780103955166
$ws.Cells.Item(6, 9).Value = 1.040142008998943
$ws.Cells.Item(6, 10).Value = 1.048616246608269
$ws.Cells.Item(6, 11).Value = 1.046468324257667
$ws.Cells.Item(6, 12).Value = 1.059968564177231
$ws.Cells.Item(6, 13).Value = 1.067134718816227
$ws.Cells.Item(6, 14).Value = 1.050105401956418

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.043938495268882
$ws.Cells.Item(7, 4).Value = 1.043847958919544
$ws.Cells.Item(7, 5).Value = 1.057292758443799
$ws.Cells.Item(7, 6).Value = 1.064431108627585
$ws.Cells.Item(7, 9).Value = 1.040074511167926
$ws.Cells.Item(7, 10).Value = 1.04843037700396
$ws.Cells.Item(7, 11).Value = 1.046311552276713
$ws.Cells.Item(7, 12).Value = 1.059723366550584
$ws.Cells.Item(7, 13).Value = 1.06684456744641
$ws.Cells.Item(7, 14).Value = 1.04991926839595

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.042696213603381
$ws.Cells.Item(8, 4).Value = 1.042943204952735
$ws.Cells.Item(8, 5).Value = 1.056022144397101
$ws.Cells.Item(8, 6).Value = 1.062974647801964
$ws.Cells.Item(8, 9).Value = 1.039790328055677
$ws.Cells.Item(8, 10).Value = 1.047653006854585
$ws.Cells.Item(8, 11).Value = 1.045655460427785
$ws.Cells.Item(8, 12).Value = 1.058698860201217
$ws.Cells.Item(8, 13).Value = 1.065632856303001
$ws.Cells.Item(8, 14).Value = 1.049140794291798

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.040510270919143
$ws.Cells.Item(9, 4).Value = 1.041351119723843
$ws.Cells.Item(9, 5).Value = 1.053789418973204
$ws.Cells.Item(9, 6).Value = 1.060416840406693
$ws.Cells.Item(9, 9).Value = 1.039281850648589
$ws.Cells.Item(9, 10).Value = 1.046281477021751
$ws.Cells.Item(9, 11).Value = 1.044496339705955
$ws.Cells.Item(9, 12).Value = 1.056895036893075
$ws.Cells.Item(9, 13).Value = 1.063501771889864
$ws.Cells.Item(9, 14).Value = 1.047767316729285

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.039055116650466
$ws.Cells.Item(10, 4).Value = 1.040291284024355
$ws.Cells.Item(10, 5).Value = 1.052305202268135
$ws.Cells.Item(10, 6).Value = 1.058717510664453
$ws.Cells.Item(10, 9).Value = 1.038937707597168
$ws.Cells.Item(10, 10).Value = 1.045366022737379
$ws.Cells.Item(10, 11).Value = 1.043721620087306
$ws.Cells.Item(10, 12).Value = 1.055693538229796
$ws.Cells.Item(10, 13).Value = 1.062083858612576
$ws.Cells.Item(10, 14).Value = 1.046850562394825

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.038425535691642
$ws.Cells.Item(11, 4).Value = 1.0398327468533
$ws.Cells.Item(11, 5).Value = 1.051663544095312
$ws.Cells.Item(11, 6).Value = 1.05798308586144
$ws.Cells.Item(11, 9).Value = 1.03878747174536
$ws.Cells.Item(11, 10).Value = 1.04496936892419
$ws.Cells.Item(11, 11).Value = 1.043385699691543
$ws.Cells.Item(11, 12).Value = 1.055173535295112
$ws.Cells.Item(11, 13).Value = 1.061470562900381
$ws.Cells.Item(11, 14).Value = 1.046453345287725

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.038191758651509
$ws.Cells.Item(12, 4).Value = 1.039662483782722
$ws.Cells.Item(12, 5).Value = 1.051425357460392
$ws.Cells.Item(12, 6).Value = 1.057710498316022
$ws.Cells.Item(12, 9).Value = 1.038731484620387
$ws.Cells.Item(12, 10).Value = 1.044821996261173
$ws.Cells.Item(12, 11).Value = 1.043260855273524
$ws.Cells.Item(12, 12).Value = 1.054980421969844
$ws.Cells.Item(12, 13).Value = 1.061242858988854
$ws.Cells.Item(12, 14).Value = 1.046305763338622

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.038241901108642
$ws.Cells.Item(13, 4).Value = 1.0396990031438
$ws.Cells.Item(13, 5).Value = 1.051476442328447
$ws.Cells.Item(13, 6).Value = 1.057768959730717
$ws.Cells.Item(13, 9).Value = 1.038743502316263
$ws.Cells.Item(13, 10).Value = 1.044853609895901
$ws.Cells.Item(13, 11).Value = 1.043287637918608
$ws.Cells.Item(13, 12).Value = 1.055021843660141
$ws.Cells.Item(13, 13).Value = 1.061291697649157
$ws.Cells.Item(13, 14).Value = 1.046337421868338

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.038406210026548
$ws.Cells.Item(14, 4).Value = 1.039818671663847
$ws.Cells.Item(14, 5).Value = 1.051643852371364
$ws.Cells.Item(14, 6).Value = 1.057960549390773
$ws.Cells.Item(14, 9).Value = 1.038782847559728
$ws.Cells.Item(14, 10).Value = 1.044957187820474
$ws.Cells.Item(14, 11).Value = 1.043375381408685
$ws.Cells.Item(14, 12).Value = 1.055157571689025
$ws.Cells.Item(14, 13).Value = 1.061451738750672
$ws.Cells.Item(14, 14).Value = 1.046441146885444

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.038507456444908
$ws.Cells.Item(15, 4).Value = 1.039892411115594
$ws.Cells.Item(15, 5).Value = 1.051747019628981
$ws.Cells.Item(15, 6).Value = 1.058078622023409
$ws.Cells.Item(15, 9).Value = 1.038807065246574
$ws.Cells.Item(15, 10).Value = 1.045021000613497
$ws.Cells.Item(15, 11).Value = 1.043429433997274
$ws.Cells.Item(15, 12).Value = 1.055241203387294
$ws.Cells.Item(15, 13).Value = 1.06155035881736
$ws.Cells.Item(15, 14).Value = 1.046505050299952

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.039096910676884
$ws.Cells.Item(16, 4).Value = 1.040321723703521
$ws.Cells.Item(16, 5).Value = 1.052347808522798
$ws.Cells.Item(16, 6).Value = 1.058766281501135
$ws.Cells.Item(16, 9).Value = 1.038947652593263
$ws.Cells.Item(16, 10).Value = 1.045392342008875
$ws.Cells.Item(16, 11).Value = 1.0437439043533
$ws.Cells.Item(16, 12).Value = 1.055728054527322
$ws.Cells.Item(16, 13).Value = 1.062124575211079
$ws.Cells.Item(16, 14).Value = 1.046876919042705

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.039466796875964
$ws.Cells.Item(17, 4).Value = 1.040591122349336
$ws.Cells.Item(17, 5).Value = 1.052724940572631
$ws.Cells.Item(17, 6).Value = 1.059198006461685
$ws.Cells.Item(17, 9).Value = 1.039035513068705
$ws.Cells.Item(17, 10).Value = 1.045625206589244
$ws.Cells.Item(17, 11).Value = 1.043941040212354
$ws.Cells.Item(17, 12).Value = 1.056033511907792
$ws.Cells.Item(17, 13).Value = 1.062484945914682
$ws.Cells.Item(17, 14).Value = 1.047110114317488

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.039682594442822
$ws.Cells.Item(18, 4).Value = 1.040748294481105
$ws.Cells.Item(18, 5).Value = 1.052945013424532
$ws.Cells.Item(18, 6).Value = 1.059449958987378
$ws.Cells.Item(18, 9).Value = 1.039086642893276
$ws.Cells.Item(18, 10).Value = 1.045761007805698
$ws.Cells.Item(18, 11).Value = 1.044055981627889
$ws.Cells.Item(18, 12).Value = 1.056211704477203
$ws.Cells.Item(18, 13).Value = 1.062695208652943
$ws.Cells.Item(18, 14).Value = 1.047246108387246

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.039756184162054
$ws.Cells.Item(19, 4).Value = 1.040801892265255
$ws.Cells.Item(19, 5).Value = 1.053020069192708
$ws.Cells.Item(19, 6).Value = 1.059535891105151
$ws.Cells.Item(19, 9).Value = 1.039104056865424
$ws.Cells.Item(19, 10).Value = 1.045807308291345
$ws.Cells.Item(19, 11).Value = 1.044095166093609
$ws.Cells.Item(19, 12).Value = 1.05627246766178
$ws.Cells.Item(19, 13).Value = 1.062766913748228
$ws.Cells.Item(19, 14).Value = 1.047292474624894

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.039427106492975
$ws.Cells.Item(20, 4).Value = 1.040562214644634
$ws.Cells.Item(20, 5).Value = 1.052684467718397
$ws.Cells.Item(20, 6).Value = 1.059151672525087
$ws.Cells.Item(20, 9).Value = 1.039026098643022
$ws.Cells.Item(20, 10).Value = 1.045600224978899
$ws.Cells.Item(20, 11).Value = 1.043919894007755
$ws.Cells.Item(20, 12).Value = 1.05600073668739
$ws.Cells.Item(20, 13).Value = 1.062446274860967
$ws.Cells.Item(20, 14).Value = 1.047085097230392

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.038357823004577
$ws.Cells.Item(21, 4).Value = 1.039783430655556
$ws.Cells.Item(21, 5).Value = 1.051594550033822
$ws.Cells.Item(21, 6).Value = 1.057904125193395
$ws.Cells.Item(21, 9).Value = 1.038771266410273
$ws.Cells.Item(21, 10).Value = 1.0449266877405
$ws.Cells.Item(21, 11).Value = 1.043349545026297
$ws.Cells.Item(21, 12).Value = 1.055117602086929
$ws.Cells.Item(21, 13).Value = 1.061404607829387
$ws.Cells.Item(21, 14).Value = 1.046410603491857

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.03768596995389
$ws.Cells.Item(22, 4).Value = 1.039294114690841
$ws.Cells.Item(22, 5).Value = 1.05091016570998
$ws.Cells.Item(22, 6).Value = 1.057120961096368
$ws.Cells.Item(22, 9).Value = 1.038609985611589
$ws.Cells.Item(22, 10).Value = 1.044502989933171
$ws.Cells.Item(22, 11).Value = 1.04299054730085
$ws.Cells.Item(22, 12).Value = 1.054562565862327
$ws.Cells.Item(22, 13).Value = 1.060750257599365
$ws.Cells.Item(22, 14).Value = 1.045986303985043

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.038042089284659
$ws.Cells.Item(23, 4).Value = 1.039553478064623
$ws.Cells.Item(23, 5).Value = 1.051272886102045
$ws.Cells.Item(23, 6).Value = 1.057536015452561
$ws.Cells.Item(23, 9).Value = 1.038695583735296
$ws.Cells.Item(23, 10).Value = 1.044727620608668
$ws.Cells.Item(23, 11).Value = 1.043180896141996
$ws.Cells.Item(23, 12).Value = 1.05485677942904
$ws.Cells.Item(23, 13).Value = 1.061097085128131
$ws.Cells.Item(23, 14).Value = 1.046211253661865

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.039445040716569
$ws.Cells.Item(24, 4).Value = 1.04057527667953
$ws.Cells.Item(24, 5).Value = 1.052702755356007
$ws.Cells.Item(24, 6).Value = 1.059172608420551
$ws.Cells.Item(24, 9).Value = 1.03903035298046
$ws.Cells.Item(24, 10).Value = 1.045611513169816
$ws.Cells.Item(24, 11).Value = 1.043929449205233
$ws.Cells.Item(24, 12).Value = 1.056015546323195
$ws.Cells.Item(24, 13).Value = 1.062463748445972
$ws.Cells.Item(24, 14).Value = 1.047096401451834

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.04107501551716
$ws.Cells.Item(25, 4).Value = 1.041762443107785
$ws.Cells.Item(25, 5).Value = 1.054365884258826
$ws.Cells.Item(25, 6).Value = 1.061077063132807
$ws.Cells.Item(25, 9).Value = 1.039414215308662
$ws.Cells.Item(25, 10).Value = 1.046636247718328
$ws.Cells.Item(25, 11).Value = 1.044796351276011
$ws.Cells.Item(25, 12).Value = 1.057361187015423
$ws.Cells.Item(25, 13).Value = 1.06624373195089
$ws.Cells.Item(25, 14).Value = 1.048122591240943
